# Update Tnf-Tnfrsf1b LR-pairs sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New numeric values (columns E:T) for rows 2-7, recomputed with new TPM data.
# Columns A:D (cluster/ligand/receptor labels) are unchanged.
$data = @{
    2 = @(3, 1, 2.913576333333333, 8.740729, 0.8649322955011439, 0.8649322955011439, 3, 1, 4.230734666666667, 12.692204, 0.3081346507358854, 0.3081346507358855, 12.32656839741289, 110.939115576716, 0.2665156107844326, 0.2665156107844326)
    3 = @(3, 1, 2.913576333333333, 8.740729, 0.8649322955011439, 0.8649322955011439, 3, 1, 8.548386000000001, 25.645158, 0.6225996527787135, 0.6225996527787135, 24.906375137798, 224.157376240182, 0.5385065468561079, 0.5385065468561079)
    4 = @(3, 1, 2.913576333333333, 8.740729, 0.8649322955011439, 0.8649322955011439, 3, 1, 0.9510283333333334, 2.853085, 0.0692656964854011, 0.0692656964854011, 2.770893644329445, 24.938042798965, 0.05991013786060349, 0.05991013786060349)
    5 = @(2, 0.6666666666666666, 0.4549836666666667, 1.364951, 0.1350677044988561, 0.1350677044988561, 3, 1, 4.230734666666667, 12.692204, 0.3081346507358854, 0.3081346507358855, 1.924915171333778, 17.324236542004, 0.0416190399514528, 0.04161903995145282)
    6 = @(2, 0.6666666666666666, 0.4549836666666667, 1.364951, 0.1350677044988561, 0.1350677044988561, 3, 1, 8.548386000000001, 25.645158, 0.6225996527787135, 0.6225996527787135, 3.889376006362, 35.004384057258, 0.08409310592260567, 0.0840931059226057)
    7 = @(2, 0.6666666666666666, 0.4549836666666667, 1.364951, 0.1350677044988561, 0.1350677044988561, 3, 1, 0.9510283333333334, 2.853085, 0.0692656964854011, 0.0692656964854011, 0.4327023582038889, 3.894321223835, 0.009355558624797609, 0.009355558624797611)
}

$columns = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $rowNum).Value = $values[$i]
    }
}

# Remove rows 8, 9 and 10 (the MuSCs sending-cluster rows are no longer present).
$ws.Range("A8:T10").EntireRow.Delete()
